$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 61.555557
$ws.Range("I4").Value = 31.75
$ws.Range("K4").Value = 31.75
$ws.Range("M4").Value = 82.25

$ws.Range("H15").Value = 593.32855
$ws.Range("I15").Value = 593.32855
$ws.Range("K15").Value = 1779.98565
$ws.Range("M15").Value = -1610.98565

$ws.Range("H17").Value = 1584039.6
$ws.Range("J17").Value = 1584039.6
$ws.Range("L17").Value = 4752118.800000001
$ws.Range("N17").Value = -4752454.800000001

$ws.Range("H32").Value = 9642.714
$ws.Range("J32").Value = 11083.167
$ws.Range("L32").Value = 11083.167
$ws.Range("N32").Value = -11735.167

$ws.Range("H70").Value = 3240.9644
$ws.Range("I70").Value = 1360
$ws.Range("J70").Value = 4285.9443
$ws.Range("K70").Value = 4080
$ws.Range("L70").Value = 12857.8329
$ws.Range("M70").Value = -3810
$ws.Range("N70").Value = -13397.8329

$ws.Range("H73").Value = 3240.9644
$ws.Range("I73").Value = 1360
$ws.Range("J73").Value = 4285.9443
$ws.Range("K73").Value = 4080
$ws.Range("L73").Value = 12857.8329
$ws.Range("M73").Value = -3144
$ws.Range("N73").Value = -14729.8329

$ws.Range("H88").Value = 914.2308
$ws.Range("I88").Value = 497.33334
$ws.Range("J88").Value = 1039.3
$ws.Range("K88").Value = 497.33334
$ws.Range("L88").Value = 1039.3
$ws.Range("M88").Value = -91.33334000000002
$ws.Range("N88").Value = -1851.3

$ws.Range("H91").Value = 914.2308
$ws.Range("I91").Value = 497.33334
$ws.Range("J91").Value = 1039.3
$ws.Range("K91").Value = 497.33334
$ws.Range("L91").Value = 1039.3
$ws.Range("M91").Value = 906.66666
$ws.Range("N91").Value = -3847.3

$ws.Range("H100").Value = 2296.1538
$ws.Range("I100").Value = 1955.5555
$ws.Range("K100").Value = 1955.5555
$ws.Range("M100").Value = -1414.5555

$ws.Range("H101").Value = 1139.091
$ws.Range("I101").Value = 1658.3334
$ws.Range("J101").Value = 516
$ws.Range("K101").Value = 4975.0002
$ws.Range("L101").Value = 1548
$ws.Range("M101").Value = -3353.0002
$ws.Range("N101").Value = -4792

$ws.Range("H135").Value = 3830.8408
$ws.Range("I135").Value = 2535.862
$ws.Range("K135").Value = 22822.758
$ws.Range("M135").Value = -20287.758

$ws.Range("H138").Value = 3419.25
$ws.Range("I138").Value = 4007.1765
$ws.Range("J138").Value = 3096.8386
$ws.Range("K138").Value = 12021.5295
$ws.Range("L138").Value = 9290.515800000001
$ws.Range("M138").Value = -6881.529500000001
$ws.Range("N138").Value = -19570.5158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8139.99
$ws.Range("I32").Value = 3581.3293
$ws.Range("J32").Value = 28907.223
$ws.Range("K32").Value = 3581.3293
$ws.Range("L32").Value = 28907.223
$ws.Range("M32").Value = -3294.3293
$ws.Range("N32").Value = -29481.223

$ws.Range("H61").Value = 11868.868
$ws.Range("I61").Value = 9392.357
$ws.Range("J61").Value = 13313.5
$ws.Range("K61").Value = 9392.357
$ws.Range("L61").Value = 13313.5
$ws.Range("M61").Value = -9180.357
$ws.Range("N61").Value = -13737.5

$ws.Range("H74").Value = 12828.054
$ws.Range("I74").Value = 1838.0834
$ws.Range("K74").Value = 1838.0834
$ws.Range("M74").Value = -964.0834

$ws.Range("H77").Value = 12828.054
$ws.Range("I77").Value = 1838.0834
$ws.Range("K77").Value = 9190.416999999999
$ws.Range("M77").Value = -4822.416999999999

$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992

$ws.Range("H119").Value = 66437.5
$ws.Range("J119").Value = 66437.5
$ws.Range("L119").Value = 66437.5
$ws.Range("N119").Value = -76113.5

$ws.Range("H122").Value = 2896.1396
$ws.Range("I122").Value = 1884.1724
$ws.Range("K122").Value = 5652.5172
$ws.Range("M122").Value = -3202.5172

$ws.Range("H132").Value = 2573971.8
$ws.Range("I132").Value = 4074.96
$ws.Range("K132").Value = 12224.88
$ws.Range("M132").Value = -9694.880000000001

$ws.Range("H136").Value = 11868.868
$ws.Range("I136").Value = 9392.357
$ws.Range("J136").Value = 13313.5
$ws.Range("K136").Value = 28177.071
$ws.Range("L136").Value = 39940.5
$ws.Range("M136").Value = -25627.071
$ws.Range("N136").Value = -45040.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 21933.334
$ws.Range("J122").Value = 21933.334
$ws.Range("L122").Value = 21933.334
$ws.Range("N122").Value = -31733.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19392.967
$ws.Range("I31").Value = 9626.23
$ws.Range("J31").Value = 26446.723
$ws.Range("K31").Value = 9626.23
$ws.Range("L31").Value = 26446.723
$ws.Range("M31").Value = -9331.23
$ws.Range("N31").Value = -27036.723

$ws.Range("H34").Value = 19392.967
$ws.Range("I34").Value = 9626.23
$ws.Range("J34").Value = 26446.723
$ws.Range("K34").Value = 9626.23
$ws.Range("L34").Value = 26446.723
$ws.Range("M34").Value = -9424.23
$ws.Range("N34").Value = -26850.723

$ws.Range("H120").Value = 18428.572
$ws.Range("J120").Value = 18428.572
$ws.Range("L120").Value = 18428.572
$ws.Range("N120").Value = -25686.572

$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -42620

$ws.Range("H122").Value = 4774.6333
$ws.Range("I122").Value = 2709.8462
$ws.Range("K122").Value = 8129.5386
$ws.Range("M122").Value = -5679.5386

$ws.Range("H125").Value = 103333.336
$ws.Range("J125").Value = 103333.336
$ws.Range("L125").Value = 103333.336
$ws.Range("N125").Value = -108253.336

$ws.Range("H132").Value = 8269.370000000001
$ws.Range("I132").Value = 3543.2666
$ws.Range("J132").Value = 14177
$ws.Range("K132").Value = 10629.7998
$ws.Range("L132").Value = 42531
$ws.Range("M132").Value = -8099.799800000001
$ws.Range("N132").Value = -47591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1391.8572
$ws.Range("I92").Value = 1799.6666
$ws.Range("K92").Value = 5398.9998
$ws.Range("M92").Value = -4150.9998

$ws.Range("H115").Value = 5037.125
$ws.Range("J115").Value = 11833.333
$ws.Range("L115").Value = 35499.999
$ws.Range("N115").Value = -37849.999

$ws.Range("H131").Value = 1482.37
$ws.Range("J131").Value = 1498.7396
$ws.Range("L131").Value = 4496.218800000001
$ws.Range("N131").Value = -14576.2188

$ws.Range("H141").Value = 5560.3
$ws.Range("I141").Value = 862.2308
$ws.Range("K141").Value = 2586.6924
$ws.Range("M141").Value = 2593.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15504

$ws.Range("H36").Value = 23557.143
$ws.Range("I36").Value = 15000
$ws.Range("J36").Value = 24983.334
$ws.Range("K36").Value = 15000
$ws.Range("L36").Value = 24983.334
$ws.Range("M36").Value = -14515
$ws.Range("N36").Value = -25953.334

$ws.Range("H102").Value = 4626.5557
$ws.Range("I102").Value = 1965.8334
$ws.Range("K102").Value = 1965.8334
$ws.Range("M102").Value = -343.8334

$ws.Range("H132").Value = 8612.416999999999
$ws.Range("I132").Value = 2849.5625
$ws.Range("J132").Value = 20138.125
$ws.Range("K132").Value = 8548.6875
$ws.Range("L132").Value = 60414.375
$ws.Range("M132").Value = -6018.6875
$ws.Range("N132").Value = -65474.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9825
$ws.Range("I7").Value = 3650.2
$ws.Range("K7").Value = 3650.2
$ws.Range("M7").Value = -3538.2

$ws.Range("H122").Value = 53772.453
$ws.Range("J122").Value = 11239.6
$ws.Range("L122").Value = 33718.8
$ws.Range("N122").Value = -38618.8

$ws.Range("H126").Value = 9825
$ws.Range("I126").Value = 3650.2
$ws.Range("K126").Value = 10950.6
$ws.Range("M126").Value = -8480.599999999999

$ws.Range("H137").Value = 83836.38
$ws.Range("J137").Value = 91285.56
$ws.Range("L137").Value = 91285.56
$ws.Range("N137").Value = -101485.56

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4318.3
$ws.Range("I62").Value = 3997
$ws.Range("J62").Value = 4639.6
$ws.Range("K62").Value = 3997
$ws.Range("L62").Value = 4639.6
$ws.Range("M62").Value = -3373
$ws.Range("N62").Value = -5887.6

$ws.Range("H65").Value = 4318.3
$ws.Range("I65").Value = 3997
$ws.Range("J65").Value = 4639.6
$ws.Range("K65").Value = 19985
$ws.Range("L65").Value = 23198
$ws.Range("M65").Value = -16865
$ws.Range("N65").Value = -29438

$ws.Range("H124").Value = 65250
$ws.Range("J124").Value = 65250
$ws.Range("L124").Value = 65250
$ws.Range("N124").Value = -75070

$ws.Range("H132").Value = 12606.359
$ws.Range("I132").Value = 5962.2085
$ws.Range("K132").Value = 17886.6255
$ws.Range("M132").Value = -15356.6255
